$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "2026-02-05 23:47:38"
$ws.Range("E3").Value = "2026-02-05 23:47:40"
$ws.Range("E4").Value = "2026-02-05 23:47:43"
$ws.Range("K4").Value = "6.5 MJ/m2"
$ws.Range("O4").Value = "12.0 °C"
$ws.Range("E5").Value = "2026-02-05 23:47:45"
$ws.Range("E6").Value = "2026-02-05 23:47:47"
$ws.Range("H6").Value = "'68%"
$ws.Range("E7").Value = "2026-02-05 23:47:50"
$ws.Range("H7").Value = "'77%"
$ws.Range("E8").Value = "2026-02-05 23:47:52"
$ws.Range("E9").Value = "2026-02-05 23:47:54"
$ws.Range("E10").Value = "2026-02-05 23:47:57"
$ws.Range("O10").Value = "7.7 °C"
$ws.Range("E11").Value = "2026-02-05 23:47:59"
$ws.Range("O11").Value = "1.2 °C"
$ws.Range("E12").Value = "2026-02-05 23:48:01"
$ws.Range("H12").Value = "'82%"
$ws.Range("E13").Value = "2026-02-05 23:48:04"
$ws.Range("O13").Value = "7.8 °C"
$ws.Range("E14").Value = "2026-02-05 23:48:06"
$ws.Range("I14").Value = "8.3 mm"
$ws.Range("E15").Value = "2026-02-05 23:48:08"
$ws.Range("H15").Value = "'77%"
$ws.Range("J15").Value = "990.3 hPa"
$ws.Range("O15").Value = "9.3 °C"
$ws.Range("E16").Value = "2026-02-05 23:48:11"
$ws.Range("H16").Value = "'97%"
$ws.Range("O16").Value = "4.1 °C"
$ws.Range("E17").Value = "2026-02-05 23:48:13"
$ws.Range("I17").Value = "9.2 mm"
$ws.Range("O17").Value = "1.2 °C"
$ws.Range("E18").Value = "2026-02-05 23:48:15"
$ws.Range("E19").Value = "2026-02-05 23:48:18"
$ws.Range("E20").Value = "2026-02-05 23:48:20"
$ws.Range("O20").Value = "-1.3 °C"
$ws.Range("E21").Value = "2026-02-05 23:48:22"
$ws.Range("E22").Value = "2026-02-05 23:48:25"
$ws.Range("E23").Value = "2026-02-05 23:48:27"
$ws.Range("E24").Value = "2026-02-05 23:48:30"
$ws.Range("H24").Value = "'74%"
$ws.Range("E25").Value = "2026-02-05 23:48:32"
$ws.Range("E26").Value = "2026-02-05 23:48:35"
$ws.Range("E27").Value = "2026-02-05 23:48:37"
$ws.Range("E28").Value = "2026-02-05 23:48:39"
$ws.Range("E29").Value = "2026-02-05 23:48:42"
$ws.Range("H29").Value = "'75%"
$ws.Range("O29").Value = "10.0 °C"
$ws.Range("E30").Value = "2026-02-05 23:48:44"
$ws.Range("O30").Value = "-1.5 °C"
$ws.Range("E31").Value = "2026-02-05 23:48:47"
$ws.Range("I31").Value = "20.6 mm"
$ws.Range("J31").Value = "994.4 hPa"
$ws.Range("E32").Value = "2026-02-05 23:48:49"
$ws.Range("O32").Value = "12.3 °C"
$ws.Range("E33").Value = "2026-02-05 23:48:51"
$ws.Range("E34").Value = "2026-02-05 23:48:54"
$ws.Range("H34").Value = "'92%"
$ws.Range("O34").Value = "4.9 °C"
$ws.Range("E35").Value = "2026-02-05 23:48:56"
$ws.Range("I35").Value = "5.6 mm"
$ws.Range("E36").Value = "2026-02-05 23:48:59"
$ws.Range("J36").Value = "992.7 hPa"
